# "further cleaning to metadata"
#
# - libraryProtocol value "E7760" -> "E7420" for every data row (K2:K27)
# - libraryProtocol (K) cells get a slightly larger (11pt), black Arial font
# - roboticLibraryPrep (L) cells become a live "=FALSE()" formula instead of
#   a hard-coded boolean literal
# - the sheet's selection is moved from L2:L27 to K2:K27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 27; $r++) {
    $kCell = $ws.Cells.Item($r, 11)   # column K - libraryProtocol
    $kCell.Value = "E7420"
    $kCell.Font.Size = 11
    $kCell.Font.Color = 0
    $kCell.WrapText = $false

    $lCell = $ws.Cells.Item($r, 12)   # column L - roboticLibraryPrep
    $lCell.Formula = "=FALSE()"
}

$ws.Range("K2:K27").Select()
